# Commit: "added range argument to read_excel"
#
# The test workbook larray/tests/data/examples.xlsx gains a new sheet,
# "pop_births_deaths", inserted right after "deaths" (and before
# "pop_missing_axis_name"). It stacks the three single-table sheets
# ("pop", "births", "deaths") on top of each other, separated by a blank
# row, so a single read_excel(..., sheet="pop_births_deaths", range=...)
# call exercising the new `range` argument can pull out each table.

$wb = $excel.ActiveWorkbook

# Grab references to the existing source sheets before we start
# inserting/moving sheets around.
$popSheet    = $wb.Worksheets.Item("pop")
$birthsSheet = $wb.Worksheets.Item("births")
$deathsSheet = $wb.Worksheets.Item("deaths")

# Insert a brand-new worksheet right after "deaths" (i.e. immediately
# before "pop_missing_axis_name"), matching the sheet order in the diff:
#   pop, births, deaths, pop_births_deaths,
#   pop_missing_axis_name, pop_missing_values, pop_narrow_format
# Worksheet.Move(target) places the sheet immediately BEFORE target, so
# we target "pop_missing_axis_name" (not "deaths") to land right after it.
$missingAxisSheet = $wb.Worksheets.Item("pop_missing_axis_name")
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "pop_births_deaths"
$newSheet.Move($missingAxisSheet)

# Re-fetch by name: after Move() the old variable binding tracks the
# worksheet that ended up at the original index, not the moved sheet.
$newSheet = $wb.Worksheets.Item("pop_births_deaths")

# Copy each source table into the new sheet, stacked with one blank
# row of separation (rows 1-7, 9-15, 17-23).
$popSheet.Range("A1:E7").Copy($newSheet.Range("A1:E7"))
$birthsSheet.Range("A1:E7").Copy($newSheet.Range("A9:E15"))
$deathsSheet.Range("A1:E7").Copy($newSheet.Range("A17:E23"))

# Restore the original active/selected sheet (the last tab,
# "pop_narrow_format", was tabSelected/active before the edit).
$lastSheet = $wb.Worksheets.Item("pop_narrow_format")
$lastSheet.Activate()
